$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update id_entrenador (column D) values for rows 2-9
# Order chosen to match the original author's shared-string insertion order
$ws.Range("D2").Value = "10"
$ws.Range("D4").Value = "13"
$ws.Range("D5").Value = "14"
$ws.Range("D6").Value = "15"
$ws.Range("D7").Value = "16"
$ws.Range("D8").Value = "17"
$ws.Range("D3").Value = "18"
$ws.Range("D9").Value = "10"

# Update the selected cell in the view
$ws.Range("D3").Select()
